# Add files via upload
# Populates the "Minute2" (G), "Second2" (H), "Rep2" (I) columns for the
# heat-2 results on the ScoreM and ScoreF sheets, and updates the
# remembered cell selection on both sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ScoreM sheet (rows 2-23): G=Minute2, H=Second2, I=Rep2
# ---------------------------------------------------------------------
$wsM = $wb.Worksheets.Item("ScoreM")

$scoreMData = @(
    @(2,  10, 0, 174),
    @(3,  10, 0, 183),
    @(4,  10, 0, 168),
    @(5,  10, 0, 233),
    @(6,  10, 0, 187),
    @(7,  10, 0, 160),
    @(8,  10, 0, 246),
    @(9,  10, 0, 202),
    @(10, 10, 0, 209),
    @(11, 10, 0, 183),
    @(12, 10, 0, 246),
    @(13, 10, 0, 170),
    @(14, 10, 0, 223),
    @(15, 10, 0, 230),
    @(16, 10, 0, 203),
    @(17, 10, 0, 216),
    @(18, 10, 0, 159),
    @(19, 10, 0, 240),
    @(20, 10, 0, 232),
    @(21, 10, 0, 233),
    @(22, 10, 0, 223),
    @(23, 10, 0, 168)
)

foreach ($row in $scoreMData) {
    $r = $row[0]
    $wsM.Cells.Item($r, 7).Value = $row[1]
    $wsM.Cells.Item($r, 8).Value = $row[2]
    $wsM.Cells.Item($r, 9).Value = $row[3]
}

# ---------------------------------------------------------------------
# ScoreF sheet (rows 2-26): G=Minute2, H=Second2, I=Rep2
# ---------------------------------------------------------------------
$wsF = $wb.Worksheets.Item("ScoreF")

$scoreFData = @(
    @(2,  10, 0, 190),
    @(3,  10, 0, 183),
    @(4,  10, 0, 210),
    @(5,  10, 0, 207),
    @(6,  10, 0, 184),
    @(7,  10, 0, 214),
    @(8,  10, 0, 217),
    @(9,  10, 0, 167),
    @(10, 10, 0, 118),
    @(11, 10, 0, 189),
    @(12, 10, 0, 192),
    @(13, 10, 0, 155),
    @(14, 10, 0, 136),
    @(15, 10, 0, 162),
    @(16, 10, 0, 164),
    @(17, 10, 0, 161),
    @(18, 10, 0, 113),
    @(19, 10, 0, 243),
    @(20, 10, 0, 205),
    @(21, 10, 0, 200),
    @(22, 10, 0, 159),
    @(23, 10, 0, 184),
    @(24, 10, 0, 178),
    @(25, 10, 0, 167),
    @(26, 10, 0, 166)
)

foreach ($row in $scoreFData) {
    $r = $row[0]
    $wsF.Cells.Item($r, 7).Value = $row[1]
    $wsF.Cells.Item($r, 8).Value = $row[2]
    $wsF.Cells.Item($r, 9).Value = $row[3]
}

# ---------------------------------------------------------------------
# Restore the active-cell selections recorded in the sheet views.
# ScoreM is updated first (it is not the active tab), then ScoreF is
# selected last so it remains the active sheet/tab, matching the
# original workbook state (ScoreF tab selected).
# ---------------------------------------------------------------------
$wsM.Range("E11").Select()
$wsF.Range("I11").Select()
